$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 154.84
$ws.Range("I15").Value = 154.84
$ws.Range("K15").Value = 464.52
$ws.Range("M15").Value = -295.52

$ws.Range("H100").Value = 2051.7
$ws.Range("I100").Value = 1098.6
$ws.Range("K100").Value = 1098.6
$ws.Range("M100").Value = -557.5999999999999

$ws.Range("H129").Value = 1033.1091
$ws.Range("I129").Value = 533.3333
$ws.Range("K129").Value = 1599.9999
$ws.Range("M129").Value = 3400.0001

$ws.Range("H137").Value = 2659.7637
$ws.Range("I137").Value = 1557.5518
$ws.Range("J137").Value = 3889.1538
$ws.Range("K137").Value = 4672.6554
$ws.Range("L137").Value = 11667.4614
$ws.Range("M137").Value = -2122.6554
$ws.Range("N137").Value = -16767.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6507.614
$ws.Range("I32").Value = 5127.2285
$ws.Range("J32").Value = 11875.777
$ws.Range("K32").Value = 5127.2285
$ws.Range("L32").Value = 11875.777
$ws.Range("M32").Value = -4840.2285
$ws.Range("N32").Value = -12449.777

$ws.Range("H96").Value = 17081.334
$ws.Range("J96").Value = 17081.334
$ws.Range("L96").Value = 17081.334
$ws.Range("N96").Value = -22573.334

$ws.Range("H97").Value = 1401.75
$ws.Range("I97").Value = 782
$ws.Range("J97").Value = 1844.4286
$ws.Range("K97").Value = 782
$ws.Range("L97").Value = 1844.4286
$ws.Range("M97").Value = -286
$ws.Range("N97").Value = -2836.4286

$ws.Range("H102").Value = 2037
$ws.Range("I102").Value = 1667.5385
$ws.Range("J102").Value = 2637.375
$ws.Range("K102").Value = 1667.5385
$ws.Range("L102").Value = 2637.375
$ws.Range("M102").Value = -45.53850000000011
$ws.Range("N102").Value = -5881.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7577704
$ws.Range("I86").Value = 8773944
$ws.Range("J86").Value = 1517.3334
$ws.Range("K86").Value = 8773944
$ws.Range("L86").Value = 1517.3334
$ws.Range("M86").Value = -8772821
$ws.Range("N86").Value = -3763.3334

$ws.Range("H89").Value = 7577704
$ws.Range("I89").Value = 8773944
$ws.Range("J89").Value = 1517.3334
$ws.Range("K89").Value = 43869720
$ws.Range("L89").Value = 7586.666999999999
$ws.Range("M89").Value = -43864104
$ws.Range("N89").Value = -18818.667

$ws.Range("H94").Value = 1281.3
$ws.Range("I94").Value = 1101.625
$ws.Range("K94").Value = 1101.625
$ws.Range("M94").Value = -650.625

$ws.Range("H134").Value = 2495.303
$ws.Range("I134").Value = 2565.125
$ws.Range("J134").Value = 2104.3
$ws.Range("K134").Value = 7695.375
$ws.Range("L134").Value = 6312.900000000001
$ws.Range("M134").Value = -5160.375
$ws.Range("N134").Value = -11382.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 267.66666
$ws.Range("I22").Value = 284.83334
$ws.Range("K22").Value = 284.83334
$ws.Range("M22").Value = 65.16665999999998

$ws.Range("H31").Value = 2349.6738
$ws.Range("I31").Value = 1893.0571
$ws.Range("K31").Value = 1893.0571
$ws.Range("M31").Value = -1598.0571

$ws.Range("H34").Value = 2349.6738
$ws.Range("I34").Value = 1893.0571
$ws.Range("K34").Value = 1893.0571
$ws.Range("M34").Value = -1691.0571

$ws.Range("H58").Value = 1936996.2
$ws.Range("I58").Value = 3368866.8
$ws.Range("J58").Value = 3971.1
$ws.Range("K58").Value = 3368866.8
$ws.Range("L58").Value = 3971.1
$ws.Range("M58").Value = -3368663.8
$ws.Range("N58").Value = -4377.1

$ws.Range("H107").Value = 1128.1333
$ws.Range("I107").Value = 1228.3846
$ws.Range("J107").Value = 476.5
$ws.Range("K107").Value = 1228.3846
$ws.Range("L107").Value = 476.5
$ws.Range("M107").Value = 691.6153999999999
$ws.Range("N107").Value = -4316.5

$ws.Range("H132").Value = 7534
$ws.Range("I132").Value = 9528.786
$ws.Range("J132").Value = 3544.4285
$ws.Range("K132").Value = 28586.358
$ws.Range("L132").Value = 10633.2855
$ws.Range("M132").Value = -26056.358
$ws.Range("N132").Value = -15693.2855

$ws.Range("H134").Value = 3275.9
$ws.Range("I134").Value = 2057.1924
$ws.Range("J134").Value = 4596.1665
$ws.Range("K134").Value = 6171.5772
$ws.Range("L134").Value = 13788.4995
$ws.Range("M134").Value = -3636.5772
$ws.Range("N134").Value = -18858.4995

$ws.Range("H136").Value = 1936996.2
$ws.Range("I136").Value = 3368866.8
$ws.Range("J136").Value = 3971.1
$ws.Range("K136").Value = 10106600.4
$ws.Range("L136").Value = 11913.3
$ws.Range("M136").Value = -10104050.4
$ws.Range("N136").Value = -17013.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1020
$ws.Range("J4").Value = 1066.6666
$ws.Range("L4").Value = 3199.9998
$ws.Range("N4").Value = -3423.9998

$ws.Range("H38").Value = 74.09524
$ws.Range("I38").Value = 23.222221
$ws.Range("K38").Value = 69.666663
$ws.Range("M38").Value = 277.333337

$ws.Range("H68").Value = 5947.2104
$ws.Range("I68").Value = 570
$ws.Range("J68").Value = 7867.643
$ws.Range("K68").Value = 1710
$ws.Range("L68").Value = 23602.929
$ws.Range("M68").Value = -899
$ws.Range("N68").Value = -25224.929

$ws.Range("H71").Value = 5947.2104
$ws.Range("I71").Value = 570
$ws.Range("J71").Value = 7867.643
$ws.Range("K71").Value = 5130
$ws.Range("L71").Value = 70808.787
$ws.Range("M71").Value = -1074
$ws.Range("N71").Value = -78920.787

$ws.Range("H80").Value = 4616.9165
$ws.Range("J80").Value = 4650.3
$ws.Range("L80").Value = 13950.9
$ws.Range("N80").Value = -15822.9

$ws.Range("H83").Value = 4616.9165
$ws.Range("J83").Value = 4650.3
$ws.Range("L83").Value = 41852.7
$ws.Range("N83").Value = -51212.7

$ws.Range("H94").Value = 3501.4443
$ws.Range("I94").Value = 999.5
$ws.Range("J94").Value = 3814.1875
$ws.Range("K94").Value = 2998.5
$ws.Range("L94").Value = 11442.5625
$ws.Range("N94").Value = -12794.5625
$ws.Range("M94").Value = -2322.5

$ws.Range("H113").Value = 745.99
$ws.Range("I113").Value = 747.67444
$ws.Range("J113").Value = 735.6429000000001
$ws.Range("K113").Value = 2243.02332
$ws.Range("L113").Value = 2206.9287
$ws.Range("M113").Value = -73.02332000000024
$ws.Range("N113").Value = -6546.9287

$ws.Range("H122").Value = 702.48
$ws.Range("J122").Value = 806.2857
$ws.Range("L122").Value = 7256.571300000001
$ws.Range("N122").Value = -12156.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1500
$ws.Range("I5").Value = 1500
$ws.Range("K5").Value = 1500
$ws.Range("M5").Value = -1388

$ws.Range("H80").Value = 8215.6
$ws.Range("I80").Value = 18450
$ws.Range("J80").Value = 3829.4285
$ws.Range("K80").Value = 18450
$ws.Range("L80").Value = 3829.4285
$ws.Range("M80").Value = -17452
$ws.Range("N80").Value = -5825.4285

$ws.Range("H83").Value = 8215.6
$ws.Range("I83").Value = 18450
$ws.Range("J83").Value = 3829.4285
$ws.Range("K83").Value = 92250
$ws.Range("L83").Value = 19147.1425
$ws.Range("M83").Value = -87258
$ws.Range("N83").Value = -29131.1425

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4223.027
$ws.Range("I132").Value = 4066.625
$ws.Range("K132").Value = 12199.875
$ws.Range("M132").Value = -9669.875

$ws.Range("H136").Value = 4767.1753
$ws.Range("I136").Value = 3505.7812
$ws.Range("J136").Value = 6381.76
$ws.Range("K136").Value = 10517.3436
$ws.Range("L136").Value = 19145.28
$ws.Range("M136").Value = -7967.3436
$ws.Range("N136").Value = -24245.28

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1726
$ws.Range("I81").Value = 1065
$ws.Range("J81").Value = 2166.6667
$ws.Range("K81").Value = 2130
$ws.Range("L81").Value = 4333.3334
$ws.Range("M81").Value = -1069
$ws.Range("N81").Value = -6455.3334

$ws.Range("H84").Value = 1726
$ws.Range("I84").Value = 1065
$ws.Range("J84").Value = 2166.6667
$ws.Range("K84").Value = 10650
$ws.Range("L84").Value = 21666.667
$ws.Range("M84").Value = -5346
$ws.Range("N84").Value = -32274.667

$ws.Range("H96").Value = 2613
$ws.Range("I96").Value = 2424.5
$ws.Range("J96").Value = 2990
$ws.Range("K96").Value = 2424.5
$ws.Range("L96").Value = 2990
$ws.Range("M96").Value = -1051.5
$ws.Range("N96").Value = -5736

$ws.Range("H122").Value = 3935
$ws.Range("I122").Value = 2362.182
$ws.Range("J122").Value = 5665.1
$ws.Range("K122").Value = 7086.545999999999
$ws.Range("L122").Value = 16995.3
$ws.Range("M122").Value = -4636.545999999999
$ws.Range("N122").Value = -21895.3

$ws.Range("H132").Value = 2374.625
$ws.Range("I132").Value = 2345.75
$ws.Range("J132").Value = 2432.375
$ws.Range("K132").Value = 7037.25
$ws.Range("L132").Value = 7297.125
$ws.Range("M132").Value = -4507.25
$ws.Range("N132").Value = -12357.125

$ws.Range("H136").Value = 4134.397
$ws.Range("I136").Value = 1793.2
$ws.Range("J136").Value = 8206.044
$ws.Range("K136").Value = 5379.6
$ws.Range("L136").Value = 24618.132
$ws.Range("M136").Value = -2829.6
$ws.Range("N136").Value = -29718.132
